$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D3").Value = "asd"
$ws.Range("E3").Value = "dsadasdasdasdasd"

$ws.Range("D1").Value = "<<ColsFit>>"
$ws.Range("E1").Value = "<<ColsFit>>"

$ws.Range("D4").Value = "dsads"
$ws.Range("E4").Value = "sdsasa"

$ws.Range("D5").Select()
